$wb = $excel.ActiveWorkbook

foreach ($ws in $wb.Worksheets) {
    # Rename column O header from "Num_Isolates" to "Present_SR"
    $ws.Range("O1").Value = "Present_SR"

    # Remove the "single_lineage" column (AO); this shifts "Phenos" (AP) left into AO
    $ws.Range("AO1").EntireColumn.Delete()
}
